# Applies the quantity/value corrections described in the commit diff for
# CryCompanywiseStockReport_1.xlsx. The workbook stores only literal values
# (no formulas), so each affected cell is written with its final target
# value directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = 9
$ws.Range("G13").Value = 15133.68
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 9690.17
$ws.Range("B15").Value = 70526.33
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 304.32
$ws.Range("F30").Value = 21
$ws.Range("G30").Value = 2151.45
$ws.Range("F33").Value = 21
$ws.Range("G33").Value = 753.0599999999999
$ws.Range("F34").Value = 38
$ws.Range("G34").Value = 1168.12
$ws.Range("B40").Value = 52249.6
$ws.Range("F43").Value = 111
$ws.Range("G43").Value = 2849.37
$ws.Range("F44").Value = 401
$ws.Range("G44").Value = 14600.41
$ws.Range("F51").Value = 74
$ws.Range("G51").Value = 1682.76
$ws.Range("F57").Value = 119
$ws.Range("G57").Value = 11131.26
$ws.Range("B72").Value = 170763.19
$ws.Range("F126").Value = 76
$ws.Range("G126").Value = 10242.52
$ws.Range("B129").Value = 67030.14
$ws.Range("B132").Value = 64196
$ws.Range("B133").Value = 65258
$ws.Range("F149").Value = 39
$ws.Range("G149").Value = 7210.71
$ws.Range("B153").Value = 18933.63
$ws.Range("F188").Value = 6
$ws.Range("G188").Value = 537.24
$ws.Range("B199").Value = 54710.34
$ws.Range("F228").Value = 256
$ws.Range("G228").Value = 4736
$ws.Range("B235").Value = 11754.35
$ws.Range("F238").Value = 19
$ws.Range("G238").Value = 2177.97
$ws.Range("B246").Value = 12802.19
$ws.Range("F250").Value = 6
$ws.Range("G250").Value = 3345
$ws.Range("B254").Value = 4687.73
$ws.Range("F274").Value = 2
$ws.Range("G274").Value = 254.42
$ws.Range("B301").Value = 95053.38
$ws.Range("B312").Value = 57802
$ws.Range("E312").Value = 162.71
$ws.Range("F312").Value = -79
$ws.Range("G312").Value = -11334.92
$ws.Range("B313").Value = 63531
$ws.Range("E313").Value = 152.53
$ws.Range("F313").Value = 23
$ws.Range("G313").Value = 3300.04
$ws.Range("B334").Value = -22948.23
$ws.Range("F351").Value = 45
$ws.Range("G351").Value = 7384.5
$ws.Range("F357").Value = 184
$ws.Range("G357").Value = 26613.76
$ws.Range("B362").Value = 69929.42999999999
$ws.Range("F366").Value = 27
$ws.Range("G366").Value = 569.97
$ws.Range("B369").Value = 58476.78
$ws.Range("F402").Value = 111
$ws.Range("G402").Value = 2828.28
$ws.Range("F406").Value = 108
$ws.Range("G406").Value = 19565.28
$ws.Range("F408").Value = 14
$ws.Range("G408").Value = 480.34
$ws.Range("F418").Value = 78
$ws.Range("G418").Value = 972.66
$ws.Range("F421").Value = 49
$ws.Range("G421").Value = 2680.3
$ws.Range("B423").Value = 154988.76
$ws.Range("F436").Value = 197
$ws.Range("G436").Value = 9117.16
$ws.Range("F438").Value = 50
$ws.Range("G438").Value = 2420.5
$ws.Range("F439").Value = 75
$ws.Range("G439").Value = 723
$ws.Range("B444").Value = 20067.92
$ws.Range("F458").Value = 47
$ws.Range("G458").Value = 12747.34
$ws.Range("F460").Value = 49
$ws.Range("G460").Value = 13865.04
$ws.Range("F461").Value = 32
$ws.Range("G461").Value = 7111.36
$ws.Range("B464").Value = 80665.22
$ws.Range("B485").Value = 53319
$ws.Range("E485").Value = 310.64
$ws.Range("F485").Value = -6
$ws.Range("G485").Value = -1643.52
$ws.Range("B486").Value = 64810
$ws.Range("E486").Value = 291.22
$ws.Range("F486").Value = 0
$ws.Range("G486").Value = 0
$ws.Range("F525").Value = 344
$ws.Range("G525").Value = 18878.72
$ws.Range("B531").Value = 106869.77
$ws.Range("F533").Value = 15
$ws.Range("G533").Value = 496.65
$ws.Range("F535").Value = 101
$ws.Range("G535").Value = 3344.11
$ws.Range("F536").Value = 8
$ws.Range("G536").Value = 345.44
$ws.Range("F537").Value = 175
$ws.Range("G537").Value = 5794.25
$ws.Range("F540").Value = 107
$ws.Range("G540").Value = 4682.32
$ws.Range("B541").Value = 18392.53
$ws.Range("F550").Value = 25
$ws.Range("G550").Value = 1547.5
$ws.Range("F552").Value = 36
$ws.Range("G552").Value = 5461.2
$ws.Range("B562").Value = 34679.16
$ws.Range("F569").Value = 9
$ws.Range("G569").Value = 1681.74
$ws.Range("B579").Value = 12605.68
$ws.Range("F611").Value = 154
$ws.Range("G611").Value = 20497.4
$ws.Range("B613").Value = 20497.4
$ws.Range("F665").Value = 27
$ws.Range("G665").Value = 1445.58
$ws.Range("B674").Value = 9586.74
$ws.Range("F680").Value = 334
$ws.Range("G680").Value = 54478.74
$ws.Range("B686").Value = 55491.29
$ws.Range("F695").Value = 23
$ws.Range("G695").Value = 2883.05
$ws.Range("B697").Value = 9022.690000000001
$ws.Range("F715").Value = 2
$ws.Range("G715").Value = 618.84
$ws.Range("B719").Value = 54872.02
$ws.Range("B724").Value = 2169989.87
$ws.Range("B725").Value = 2169989.87
